# edit.ps1 - COM-interop script applying the commit's changes:
#   1. Slide 6's table switches from the custom "Table_0" style to the
#      built-in "No Style, Table Grid" style.
#   2. The slide master's theme (ppt/theme/theme1.xml) swaps its color
#      scheme from "Integral" to the stock "Office" palette (the theme's
#      font scheme / format scheme are identical between the two themes,
#      so only the 12 theme colors need to move).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Table style id change on slide 6 (the table is shape 2).
# ---------------------------------------------------------------------
$slide = $p.Slides.Item(6)
$tableShape = $slide.Shapes.Item(2)
$tableShape.Table.ApplyStyle("{F9F30EC8-4CCC-4F7E-9295-2E02A6B029EB}")

# ---------------------------------------------------------------------
# 2. Swap the slide master theme's colour scheme (Integral -> Office).
# ---------------------------------------------------------------------
function Convert-HexToBgrLong($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + $g * 256 + $b * 65536
}

# Order matches ThemeColorScheme.Item(1..12): dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink.
$officeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$colorScheme = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $colorScheme.Item($i).RGB = Convert-HexToBgrLong $officeColors[$i - 1]
}
